$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for A2:B12 (Transaccion, Cuotas)
$data = @(
    @(828271183, 1),
    @(668066389, 1),
    @(669436101, 1),
    @(550712337, 1),
    @(439469821, 1),
    @(488102342, 1),
    @(985252614, 1),
    @(332987301, 1),
    @(638201176, 1),
    @(592586767, 1),
    @(961249559, 3)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Remove rows 13 and 14 which no longer exist (shrinks used range / dimension to A1:B12)
$ws.Range("A13:B14").Delete()
